# check null and undefined for Model value
#
# Previously, rows with missing/null/undefined numeric metrics were written
# out as the literal text "-" in the impression/clicks/totalUser/cost
# columns. After fixing the null/undefined check, those placeholders are
# now written as the number 0 instead. This script rewrites every cell that
# currently holds the text "-" (in a column that is otherwise numeric for
# its row's src) to the numeric value 0, matching the corrected behavior.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToFix = @(
    "D44","E44","H44",
    "D131","E131","H131",
    "G162",
    "D189","E189","G189",
    "D190","E190","G190",
    "D191","E191","G191",
    "D192","E192","G192",
    "D193","E193","G193",
    "G194",
    "G196",
    "E197","G197",
    "D198","E198","G198",
    "E199","G199",
    "G200",
    "G201",
    "E202","G202",
    "E203","G203",
    "E204","G204",
    "E205","G205",
    "E206","G206",
    "E207","G207",
    "E208","G208",
    "G209",
    "E210","G210",
    "E211","G211",
    "E212","G212",
    "G213",
    "G214",
    "G215",
    "E216","G216",
    "E217","G217",
    "G218",
    "E219","G219",
    "E220","G220",
    "G221",
    "E222","G222",
    "E223","G223",
    "G224",
    "G225",
    "E226","G226",
    "G229",
    "G231",
    "G233",
    "E235","G235",
    "E236","G236",
    "G237",
    "G238",
    "E239","G239",
    "G240",
    "E242","G242",
    "E244","G244",
    "G245",
    "G246",
    "G247",
    "G248",
    "G249",
    "G250",
    "G251",
    "G252",
    "G253",
    "G254",
    "G255",
    "G256",
    "G257",
    "G258",
    "G260",
    "G261",
    "G262",
    "G264",
    "G265"
)

foreach ($addr in $cellsToFix) {
    $ws.Range($addr).Value = 0
}
